# Second commit: add a header row (name / age / address) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A1").Value = "name "
$ws.Range("B1").Value = "age"
$ws.Range("C1").Value = "address"

# After typing the last header cell, Excel's selection moved to the next row.
$ws.Range("A2:XFD2").Select()
